# This deck is an Okapi-filter OOXML regression fixture. Every run (and
# field) in it carries a literal "merge placeholder" payload of the shape
#
#   {gN:-ERR:REF-NOT-FOUND-}<visible text>{/gM:&lt;/a:t&gt;&lt;/a:r&gt;}
#
# inside its <a:t>. The "-ERR:REF-NOT-FOUND-" token is a stand-in the
# filter emits when it cannot resolve the {gN: ... } opening-tag
# reference; fixing it means substituting, in place of that token, the
# XML-escaped text of the run's (or field's) own opening tag sequence up
# to and including "<a:t>" - i.e. what {gN: ...} should have pointed to.
#
# Because PowerPoint's automation surface here has no first-class "Runs"
# collection that addresses an individual run's XML, we locate each
# occurrence of the marker text with TextRange.Text/IndexOf and rewrite
# just that sub-string via TextRange2.Characters(start,len).Text - this
# is the scripted equivalent of selecting the broken token and retyping
# the correct replacement, run by run, left to right.

function Fix-Markers($Shape, $Replacements) {
    $marker = "-ERR:REF-NOT-FOUND-"
    $tf = $Shape.TextFrame
    $tf2 = $Shape.TextFrame2

    for ($k = 0; $k -lt $Replacements.Length; $k++) {
        $text = $tf.TextRange.Text
        $idx = $text.IndexOf($marker)
        if ($idx -lt 0) {
            break
        }
        $sub = $tf2.TextRange.Characters($idx + 1, $marker.Length)
        $sub.Text = $Replacements[$k]
    }
}

# Builds the alternating (no-err / err="1") replacement list used by the
# slide's "Content Placeholder 2" shape: three identical 31-entry runs,
# one per visual line (the lines are separated by <a:br/>, not by
# paragraph breaks, and the alternation restarts at the top of each
# line).
function Build-AlternatingList($BlockSize, $BlockCount, $EvenRepl, $OddRepl) {
    $list = @()
    for ($b = 0; $b -lt $BlockCount; $b++) {
        for ($i = 0; $i -lt $BlockSize; $i++) {
            if ($i % 2 -eq 0) {
                $list += $EvenRepl
            } else {
                $list += $OddRepl
            }
        }
    }
    return $list
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Slide 1 shapes ----------------------------------------------------

# Shape "Title 1": single run, lang="fr-FR" dirty="0" smtClean="0"
$titleShape = $s.Shapes.Item(1)
$titleRepl = @(
    '<a:r><a:rPr lang="fr-FR" dirty="0" smtClean="0"/><a:t>'
)
Fix-Markers $titleShape $titleRepl

# Shape "Content Placeholder 2": 93 runs alternating between
# lang="fr-FR" dirty="0" smtClean="0"   (no err attr)
# and
# lang="fr-FR" dirty="0" err="1" smtClean="0"
$contentShape = $s.Shapes.Item(2)
$evenRepl = '<a:r><a:rPr lang="fr-FR" dirty="0" smtClean="0"/><a:t>'
$oddRepl  = '<a:r><a:rPr lang="fr-FR" dirty="0" err="1" smtClean="0"/><a:t>'
$contentRepl = Build-AlternatingList 31 3 $evenRepl $oddRepl
Fix-Markers $contentShape $contentRepl

# ---- Slide master shapes ------------------------------------------------

$master = $s.Master

# "Title Placeholder 1": single run, lang="en-US" smtClean="0"
$masterTitleRepl = @(
    '<a:r><a:rPr lang="en-US" smtClean="0"/><a:t>'
)
Fix-Markers $master.Shapes.Item(1) $masterTitleRepl

# "Text Placeholder 2": 5 runs (outline levels 0-4), all
# lang="en-US" smtClean="0"
$masterTextRepl = @(
    '<a:r><a:rPr lang="en-US" smtClean="0"/><a:t>',
    '<a:r><a:rPr lang="en-US" smtClean="0"/><a:t>',
    '<a:r><a:rPr lang="en-US" smtClean="0"/><a:t>',
    '<a:r><a:rPr lang="en-US" smtClean="0"/><a:t>',
    '<a:r><a:rPr lang="en-US" smtClean="0"/><a:t>'
)
Fix-Markers $master.Shapes.Item(2) $masterTextRepl

# "Date Placeholder 3": a:fld datetimeFigureOut
$masterDateRepl = @(
    '<a:fld id="{1D8BD707-D9CF-40AE-B4C6-C98DA3205C09}" type="datetimeFigureOut"><a:rPr lang="en-US" smtClean="0"/><a:pPr/><a:t>'
)
Fix-Markers $master.Shapes.Item(3) $masterDateRepl

# "Footer Placeholder 4" has no marker text - nothing to do.

# "Slide Number Placeholder 5": a:fld slidenum
$masterSlideNumRepl = @(
    '<a:fld id="{B6F15528-21DE-4FAA-801E-634DDDAF4B2B}" type="slidenum"><a:rPr lang="en-US" smtClean="0"/><a:pPr/><a:t>'
)
Fix-Markers $master.Shapes.Item(5) $masterSlideNumRepl
